$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.461.28"
$ws.Range("E2").Value = "  +1.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.390.17"
$ws.Range("E3").Value = "  +4.10%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.44"
$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.70"
$ws.Range("E6").Value = "  +8.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.390.14"
$ws.Range("E8").Value = "  +4.22%  "

$ws.Range("E9").Value = "  +1.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.57"
$ws.Range("E10").Value = "  +4.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("E11").Value = "  +8.15%  "

$ws.Range("E12").Value = "  +6.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.963.83"
$ws.Range("E13").Value = "  +3.54%  "

$ws.Range("E14").Value = "  +1.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").Value = "  +7.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.398.78"
$ws.Range("E16").Value = "  +3.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.38"
$ws.Range("E17").Value = "  +5.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.555.76"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.06"
$ws.Range("E19").Value = "  +6.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.89"
$ws.Range("E20").Value = "  +5.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.29"
$ws.Range("E21").Value = "  +3.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.58"
$ws.Range("E22").Value = "  +11.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.572"
$ws.Range("E23").Value = "  +3.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.529.43"
$ws.Range("E24").Value = "  +4.01%  "

$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  +18.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.04"
$ws.Range("E27").Value = "  +2.83%  "

$ws.Range("E28").Value = "  +14.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.71"
$ws.Range("E29").Value = "  +7.20%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.31"
$ws.Range("E31").Value = "  +7.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.160"
$ws.Range("E32").Value = "  +7.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.15"
$ws.Range("E33").Value = "  +3.45%  "

$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.422.91"
$ws.Range("E35").Value = "  +3.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.51"
$ws.Range("E36").Value = "  +4.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.50"
$ws.Range("E37").Value = "  +5.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.98"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("E39").Value = "  +5.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "161.80"
$ws.Range("E40").Value = "  +2.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0792"
$ws.Range("E41").Value = "  +6.04%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.73"
$ws.Range("E43").Value = "  +12.28%  "

$ws.Range("E44").Value = "  +8.58%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.773"
$ws.Range("E45").Value = "  +5.58%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.46"
$ws.Range("E46").Value = "  +3.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.22"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.18"
$ws.Range("E48").Value = "  +7.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.97"
$ws.Range("E49").Value = "  +5.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.96"
$ws.Range("E50").Value = "  +8.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.375.12"
$ws.Range("E51").Value = "  +10.92%  "
